# Update "想去人数" (want-to-go count) figures on the 展览 and 全部类型 sheets
# to reflect newly generated output data.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 3801
$ws1.Range("F5").Value = 449
$ws1.Range("F8").Value = 184
$ws1.Range("F10").Value = 95
$ws1.Range("F11").Value = 1414
$ws1.Range("F13").Value = 2385
$ws1.Range("F14").Value = 169

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 3801
$ws4.Range("F5").Value = 449
$ws4.Range("F9").Value = 184
$ws4.Range("F11").Value = 95
$ws4.Range("F14").Value = 1414
$ws4.Range("F16").Value = 2385
$ws4.Range("F17").Value = 169

$wb.Save()
